# Applies the "Biodigester O&M Cost" update to the Results_Summary workbook.
#
# Summary of the change:
#  - A new unique string "Biodigester O&M Cost" is introduced.
#  - Sheet "Project Total Costs" (sheet1): values B2:B7 updated (re-run simulation numbers).
#  - Sheet "Components Capacity and Cost" (sheet2): values B2:B12 updated.
#  - Sheet "Yearly Costs Info" (sheet3): a new column is inserted for
#    "Biodigester O&M Cost" (between "Renewable O&M Cost" and "Total O&M Cost"),
#    pushing the old "Total O&M Cost" / "Fuel Cost" / "Battery Replacement Cost" /
#    "Lost Load Cost" columns one place to the right, and all numeric values are
#    refreshed with newly computed results.
#  - Sheet "Yearly Energy Averages" (sheet4): values B2:E6 updated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Project Total Costs
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Project Total Costs")
$ws1.Cells.Item(2,2).Value = 580154.753044
$ws1.Cells.Item(3,2).Value = 93113.9803754
$ws1.Cells.Item(4,2).Value = 148508.240061
$ws1.Cells.Item(5,2).Value = 600000.111121
$ws1.Cells.Item(6,2).Value = 112959.338452
$ws1.Cells.Item(7,2).Value = 0.8180249723534793

# ---------------------------------------------------------------------------
# Sheet 2: Components Capacity and Cost
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Components Capacity and Cost")
$ws2.Cells.Item(2,2).Value = 406.638746988
$ws2.Cells.Item(3,2).Value = 7.95641204699
$ws2.Cells.Item(4,2).Value = 29.9732158777
$ws2.Cells.Item(5,2).Value = 146.925420661
$ws2.Cells.Item(6,2).Value = [double]"2.89581149031e-06"
$ws2.Cells.Item(7,2).Value = 219584.92337352
$ws2.Cells.Item(8,2).Value = 1591.282409398
$ws2.Cells.Item(9,2).Value = 7493.303969425
$ws2.Cells.Item(10,2).Value = 293850.841322
$ws2.Cells.Item(11,2).Value = 0.008687434470930001
$ws2.Cells.Item(12,2).Value = 522520.3597617776

# ---------------------------------------------------------------------------
# Sheet 3: Yearly Costs Info
# Insert a new column E for "Biodigester O&M Cost"; this shifts the previous
# E,F,G,H columns (Total O&M Cost, Fuel Cost, Battery Replacement Cost,
# Lost Load Cost) one column to the right, into F,G,H,I.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Yearly Costs Info")
$ws3.Columns.Item(5).Insert()

$ws3.Cells.Item(1,2).Value = "Battery O&M Cost"
$ws3.Cells.Item(1,3).Value = "Generator O&M Cost"
$ws3.Cells.Item(1,4).Value = "Renewable O&M Cost"
$ws3.Cells.Item(1,5).Value = "Biodigester O&M Cost"
$ws3.Cells.Item(1,6).Value = "Total O&M Cost"
$ws3.Cells.Item(1,7).Value = "Fuel Cost"
$ws3.Cells.Item(1,8).Value = "Battery Replacement Cost"
$ws3.Cells.Item(1,9).Value = "Lost Load Cost"
$ws3.Range("I1").Style = $ws3.Range("H1").Style

# Row 2 (Year 1)
$ws3.Cells.Item(2,2).Value = 4391.698467470401
$ws3.Cells.Item(2,3).Value = 408.8063870470349
$ws3.Cells.Item(2,4).Value = 5877.01700018869
$ws3.Cells.Item(2,5).Value = 3486.5892000072
$ws3.Cells.Item(2,6).Value = 14164.11105471333
$ws3.Cells.Item(2,7).Value = 75.72664166297506
$ws3.Cells.Item(2,8).Value = 5850.092220579984
$ws3.Cells.Item(2,9).Value = 0

# Row 3 (Year 2)
$ws3.Cells.Item(3,2).Value = 4391.698467470401
$ws3.Cells.Item(3,3).Value = 408.8063870470349
$ws3.Cells.Item(3,4).Value = 5877.01700018869
$ws3.Cells.Item(3,5).Value = 3486.5892000072
$ws3.Cells.Item(3,6).Value = 14164.11105471333
$ws3.Cells.Item(3,7).Value = 960.6693361899218
$ws3.Cells.Item(3,8).Value = 6936.58399086376
$ws3.Cells.Item(3,9).Value = 0

# Row 4 (Year 3)
$ws3.Cells.Item(4,2).Value = 4391.698467470401
$ws3.Cells.Item(4,3).Value = 408.8063870470349
$ws3.Cells.Item(4,4).Value = 5877.01700018869
$ws3.Cells.Item(4,5).Value = 3486.5892000072
$ws3.Cells.Item(4,6).Value = 14164.11105471333
$ws3.Cells.Item(4,7).Value = 983.9071927690284
$ws3.Cells.Item(4,8).Value = 6933.451650374377
$ws3.Cells.Item(4,9).Value = 0

# Row 5 (Year 4)
$ws3.Cells.Item(5,2).Value = 4391.698467470401
$ws3.Cells.Item(5,3).Value = 408.8063870470349
$ws3.Cells.Item(5,4).Value = 5877.01700018869
$ws3.Cells.Item(5,5).Value = 3486.5892000072
$ws3.Cells.Item(5,6).Value = 14164.11105471333
$ws3.Cells.Item(5,7).Value = 13329.86551201737
$ws3.Cells.Item(5,8).Value = 6621.112361777014
$ws3.Cells.Item(5,9).Value = 0

# Row 6 (Year 5)
$ws3.Cells.Item(6,2).Value = 4391.698467470401
$ws3.Cells.Item(6,3).Value = 408.8063870470349
$ws3.Cells.Item(6,4).Value = 5877.01700018869
$ws3.Cells.Item(6,5).Value = 3486.5892000072
$ws3.Cells.Item(6,6).Value = 14164.11105471333
$ws3.Cells.Item(6,7).Value = 29948.96724055574
$ws3.Cells.Item(6,8).Value = 6047.308640889112
$ws3.Cells.Item(6,9).Value = 0

# ---------------------------------------------------------------------------
# Sheet 4: Yearly Energy Averages
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Yearly Energy Averages")

$ws4.Cells.Item(2,2).Value = 0.0003660717592900819
$ws4.Cells.Item(2,3).Value = 0.2340607422158436
$ws4.Cells.Item(2,4).Value = 106.4830380499204
$ws4.Cells.Item(2,5).Value = 41.42708421087878

$ws4.Cells.Item(3,2).Value = 0.0004840680156889537
$ws4.Cells.Item(3,3).Value = 2.338156506860181
$ws4.Cells.Item(3,4).Value = 104.0765624558221
$ws4.Cells.Item(3,5).Value = 26.80640825805316

$ws4.Cells.Item(4,2).Value = 0.0005113701993705007
$ws4.Cells.Item(4,3).Value = 2.39469588531572
$ws4.Cells.Item(4,4).Value = 104.0765592830033
$ws4.Cells.Item(4,5).Value = 26.79598487055777

$ws4.Cells.Item(5,2).Value = 0.04308053126475316
$ws4.Cells.Item(5,3).Value = 22.07834912608757
$ws4.Cells.Item(5,4).Value = 82.07347285167728
$ws4.Cells.Item(5,5).Value = 12.93512144022593

$ws4.Cells.Item(6,2).Value = 1.160869903359796
$ws4.Cells.Item(6,3).Value = 35.06415481030373
$ws4.Cells.Item(6,4).Value = 66.51423669167458
$ws4.Cells.Item(6,5).Value = 4.670275651136829

$wb.Save()
